$d = $word.ActiveDocument

# 1. "Feature extraction (...)" line: collapse the verbose "sent freq" wording
#    down to the shorter feature list, and add "BOW". Done as several
#    smaller replacements (rather than one giant one) so each swallows a
#    "sent freq" phrase together with its neighbouring spell-check markup,
#    instead of leaving orphaned <w:proofErr/> elements behind.
$d.Content.Find.Execute(
    "POS sent frequency, NER group sent freq,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NER, BOW,", 2)

$d.Content.Find.Execute(
    "Bigram sent freq, Trigra",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bigram, Trigra", 2)

$d.Content.Find.Execute(
    "m sent freq, Quad",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "m, Quad", 2)

$d.Content.Find.Execute(
    "gram sent freq) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "gram) ", 2)

# 2. "Modeling (...)" line: drop the "Machine Learning" alternative, keep
#    "logistic model" as the sole option.
$d.Content.Find.Execute(
    "Machine Learning or logistic model",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "logistic model", 2)
